# Move the "Q&A" slide (currently slide 23) so it comes after the
# "Continuous machine learning." slide (currently slide 27), i.e. to
# slide position 27. Slides 24-27 shift up by one to fill the gap.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(23)
$s.MoveTo(27)
